$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates stored as plain text (e.g. "27/07/2025"), not real
# Excel dates, so force text formatting before assigning the value to stop
# Excel from auto-converting the string into a date serial number. Clear
# the formatting afterwards so the cell doesn't end up with a stray style.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "07/08/2025"
$ws.Range("A21").ClearFormats()

$ws.Range("B21").Value = "Godoy Cruz"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = "Gimnasia L.P."
$ws.Range("F21").Value = "L"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 2
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1.03
$ws.Range("L21").Value = 1.4
$ws.Range("M21").Value = 20
$ws.Range("N21").Value = 5
$ws.Range("O21").Value = 4
$ws.Range("P21").Value = 3
